$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B69 changes from text "2" to a real number 2
$ws.Range("B69").Value = 2

# New row 70
$ws.Range("A70").Value = "Ying Tang"

# B70 must stay stored as text "4" (not a number), matching the target diff's
# inlineStr cell, so force a text number-format before assigning the value.
$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = "4"

$ws.Range("C70").Value = "could be a game changer"
$ws.Range("D70").Value = "DIS"
$ws.Range("E70").Value = "RES"
$ws.Range("F70").Value = "a3c87a5d-b7d4-4eb1-9136-458357f6153b"
$ws.Range("G70").Value = "IrVvIL2BaXrg4_annotated.xlsx"
$ws.Range("H70").Value = "If this is confirmed over benchmark dataset this could be a game changer."
